# Copy in EU-2024-develop branch
# Updates the "About" sheet's source citation (year + linked article) and
# refreshes the BNEF battery-cost-decline assumption on "PDiBCpDoC".

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsData  = $wb.Worksheets.Item("PDiBCpDoC")

# --- "About" sheet -------------------------------------------------------
$wsAbout.Activate()

# Source year: 2019 -> 2024
$wsAbout.Range("B4").Value = 2024

# Source citation swaps from the old blog post to the new EV Outlook report
$wsAbout.Range("B6").Value = "https://about.bnef.com/electric-vehicle-outlook/"
$wsAbout.Range("B5").Value = "Electric Vehicle Outlook 2024"

# Drop the stray formatted-but-empty cell that used to live at D14
$wsAbout.Range("D14").Clear()

# Leave the selection where the edit happened
$wsAbout.Range("B6").Select()

# --- "PDiBCpDoC" sheet ----------------------------------------------------
$wsData.Activate()

# Updated decline-per-doubling assumption: 0.18 -> 0.17
$wsData.Range("B2").Value = 0.17

$wsData.Range("B3").Select()

# Restore "About" as the active sheet/tab
$wsAbout.Activate()
